$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 16.02.2022 12:30"

# Row 3 (Tesco): swap Cena (B3) and Old Cena (C3) values
$ws.Range("B3").Value = 36.5
$ws.Range("C3").Value = 35.9

# Update Delta Cena (D3) and Old Datum (E3) text values
# Force text format so the leading "+" sign is preserved as a string, not parsed as a number,
# then restore the default "Normal" style so no extra cell formatting is left behind.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "+0.6"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "2022-02-16 12:31:55"
